$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: F6 text changes; D6 and E6 are new cells added
$ws.Range("F6").Value = "Modifications pour relance facture et pdf relance facture.Ajout de la Carte"
$ws.Range("D6").Value = 0.89236111111111116
$ws.Range("D6").NumberFormat = $ws.Range("B6").NumberFormat
$ws.Range("E6").Value = "9h46"

# Row 7: B7 and F7 are new cells added
$ws.Range("B7").Value = 0.35902777777777778
$ws.Range("B7").NumberFormat = $ws.Range("B6").NumberFormat
$ws.Range("F7").Value = "Gestion déconnexion. Gestion 404. Ajout Prix & Options."

# Update selection to F7
$ws.Range("F7").Select()
